$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country rows whose shared-string order changed (value/meaning swap) ---
# Chile(31)/Chequia(32)/India(33)/Dinamarca(34): Chequia and India swap positions.
# Row 28 previously showed Chequia's numbers; it now shows India (with updated stats).
# Row 29 previously showed India's numbers; it now shows Chequia (its stats unchanged).
$ws.Range("A28").Value = "India"
$ws.Range("A29").Value = "Chequia"

# Armenia(72)/Hungria(73)/Barein(74)/Crucero(75): Hungria and Barein swap positions.
# Row 69 previously showed Hungria's numbers; it now shows Barein (with updated stats).
# Row 70 previously showed Barein's numbers; it now shows Hungria (its stats unchanged).
$ws.Range("A69").Value = "Barein"
$ws.Range("A70").Value = "Hungria"

# --- Update the "last updated" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 6 de Abril de 2020 a las 19:52"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 355834
$ws.Range("C4").Value = 19161
$ws.Range("E4").Value = 326129
$ws.Range("G4").Value = 842
$ws.Range("H4").Value = 10458

# --- Row 7 ---
$ws.Range("B7").Value = 101214
$ws.Range("C7").Value = 1091
$ws.Range("E7").Value = 70902

# --- Row 16 ---
$ws.Range("B16").Value = 16498
$ws.Range("C16").Value = 986
$ws.Range("E16").Value = 13049
$ws.Range("G16").Value = 41
$ws.Range("H16").Value = 321

# --- Row 25 ---
$ws.Range("B25").Value = 5762
$ws.Range("C25").Value = 75
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 76

# --- Row 28 (now India) updated stats ---
$ws.Range("B28").Value = 4778
$ws.Range("C28").Value = 489
$ws.Range("D28").Value = 346
$ws.Range("E28").Value = 4303
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = 129

# --- Row 29 (now Chequia) stats (unchanged values, relocated from old row 28) ---
$ws.Range("B29").Value = 4735
$ws.Range("C29").Value = 148
$ws.Range("D29").Value = 121
$ws.Range("E29").Value = 4536
$ws.Range("F29").Value = 84
$ws.Range("G29").Value = 11
$ws.Range("H29").Value = 78

# --- Row 60 ---
$ws.Range("B60").Value = 1120
$ws.Range("C60").Value = 99
$ws.Range("D60").Value = 81
$ws.Range("E60").Value = 959
$ws.Range("G60").Value = 10
$ws.Range("H60").Value = 80

# --- Row 69 (now Barein) updated stats ---
$ws.Range("B69").Value = 756
$ws.Range("C69").Value = 56
$ws.Range("D69").Value = 458
$ws.Range("E69").Value = 294
$ws.Range("F69").Value = 4
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 4

# --- Row 70 (now Hungria) stats (unchanged values, relocated from old row 69) ---
$ws.Range("B70").Value = 744
$ws.Range("C70").Value = 11
$ws.Range("D70").Value = 67
$ws.Range("E70").Value = 639
$ws.Range("F70").Value = 17
$ws.Range("G70").Value = 4
$ws.Range("H70").Value = 38
